# Daily attendance processing - rotate "Recorded By" (column G) log entries
# by moving the most-recently-appended name to the front of the list,
# for rows whose value matches one of the known trigger patterns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$triggerValues = @(
    "System, dnasr281@gmail.com",
    "backup@backdoor.com, System, system"
)

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val) { continue }

    if ($triggerValues -contains $val) {
        $parts = $val -split ", "
        $n = $parts.Length
        $lastItem = $parts[$n - 1]
        $rest = $parts[0..($n - 2)]
        $rotated = @($lastItem) + $rest
        $newVal = $rotated -join ", "
        $cell.Value = $newVal
    }
}
